# Word COM-interop script implementing the "Hide comments button added" edit.
#
# 1) The whole "Iznad komentara dodati dugme..." paragraph (describing the
#    new Hide/Show comments button) is highlighted red to flag it as new.
# 2) The old "Sta se desava sa komentarima izbrisanih postova? I oni treba
#    da nestanu..." bullet is repurposed: its text becomes
#    "Na single-post strani naslov bloga ne treba da bude link", and the
#    original sentence is moved into a brand-new bullet right after it,
#    highlighted red (same treatment as point 1 - a freshly added task).

$d = $word.ActiveDocument

$wdRed = 6

# --- 1) Highlight the whole "Iznad komentara..." paragraph red ---------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.StartsWith("Iznad")) {
        $p.Range.HighlightColorIndex = $wdRed
        break
    }
}

# --- 2) Rework the "Sta se desava..." bullet ----------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.StartsWith("Sta se desava")) {

        $oldText = $p.Range.Text
        # Strip the trailing paragraph mark so we keep it (and the
        # _GoBack bookmark sitting right before it) in place.
        $oldTextNoMark = $oldText.Substring(0, $oldText.Length - 1)

        # Insert a brand-new bullet right after this one, matching the
        # same list formatting, then fill it with the original sentence
        # and highlight it red.
        $p.Range.InsertParagraphAfter()

        $newP = $d.Paragraphs($i + 1)
        $newP.Range.Text = $oldTextNoMark
        $newP.Range.HighlightColorIndex = $wdRed

        # Now replace the text of the original bullet with the new
        # "single-post strani naslov bloga..." sentence (paragraph mark
        # and bookmark are preserved because we stop one char short of
        # the end of the range).
        $body = $d.Range($p.Range.Start, $p.Range.End - 1)
        $body.Text = "Na single-post strani naslov bloga ne treba da bude link"

        break
    }
}
